$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume refresh -- updates the Price (D) and
# Volume(1h) (E) columns for each coin row to the latest scraped values.
# Cells already hold plain text (e.g. "29.772.11", "  -0.54%  "), so we
# force the Text number format before writing to stop Excel from
# auto-coercing numeric-looking strings (like "1.000" or "24.43") into
# actual numbers, then restore General to match the original styling.
$updates = @(
    @{ Cell = 'D2'; Value = '29.795.34' }
    @{ Cell = 'E2'; Value = '  -0.41%  ' }
    @{ Cell = 'D3'; Value = '1.870.40' }
    @{ Cell = 'E3'; Value = '  -0.22%  ' }
    @{ Cell = 'D4'; Value = '1.000' }
    @{ Cell = 'E4'; Value = '  -0.08%  ' }
    @{ Cell = 'D5'; Value = '0.7279' }
    @{ Cell = 'E5'; Value = '  -1.74%  ' }
    @{ Cell = 'D6'; Value = '241.23' }
    @{ Cell = 'E6'; Value = '  -0.44%  ' }
    @{ Cell = 'D7'; Value = '1.000' }
    @{ Cell = 'E7'; Value = '  -0.04%  ' }
    @{ Cell = 'D8'; Value = '0.3131' }
    @{ Cell = 'E8'; Value = '  -0.67%  ' }
    @{ Cell = 'D9'; Value = '0.07140' }
    @{ Cell = 'E9'; Value = '  -0.59%  ' }
    @{ Cell = 'D10'; Value = '24.43' }
    @{ Cell = 'E10'; Value = '  -1.19%  ' }
    @{ Cell = 'D11'; Value = '0.08131' }
    @{ Cell = 'E11'; Value = '  -4.03%  ' }
    @{ Cell = 'D12'; Value = '0.7426' }
    @{ Cell = 'E12'; Value = '  -1.37%  ' }
    @{ Cell = 'D13'; Value = '1.880.60' }
    @{ Cell = 'E13'; Value = '  +0.44%  ' }
    @{ Cell = 'D14'; Value = '5.340' }
    @{ Cell = 'E14'; Value = '  -1.02%  ' }
    @{ Cell = 'D15'; Value = '92.38' }
    @{ Cell = 'E15'; Value = '  -0.24%  ' }
    @{ Cell = 'D16'; Value = '29.787.25' }
    @{ Cell = 'E16'; Value = '  -0.52%  ' }
    @{ Cell = 'D18'; Value = '247.48' }
    @{ Cell = 'E18'; Value = '  +1.72%  ' }
    @{ Cell = 'D19'; Value = '13.37' }
    @{ Cell = 'E19'; Value = '  -1.77%  ' }
    @{ Cell = 'D20'; Value = '0.000007802' }
    @{ Cell = 'E20'; Value = '  -0.29%  ' }
    @{ Cell = 'D21'; Value = '0.9990' }
    @{ Cell = 'E21'; Value = '  -0.12%  ' }
    @{ Cell = 'D22'; Value = '2.117.70' }
    @{ Cell = 'E22'; Value = '  -0.92%  ' }
    @{ Cell = 'D23'; Value = '1.001' }
    @{ Cell = 'E23'; Value = '  +0.26%  ' }
    @{ Cell = 'D24'; Value = '7.742' }
    @{ Cell = 'E24'; Value = '  -3.13%  ' }
    @{ Cell = 'D25'; Value = '0.1533' }
    @{ Cell = 'E25'; Value = '  -1.52%  ' }
    @{ Cell = 'D26'; Value = '9.204' }
    @{ Cell = 'E26'; Value = '  -1.12%  ' }
    @{ Cell = 'E27'; Value = '  -1.03%  ' }
    @{ Cell = 'D28'; Value = '18.53' }
    @{ Cell = 'E28'; Value = '  -0.49%  ' }
    @{ Cell = 'D29'; Value = '2.006' }
    @{ Cell = 'E29'; Value = '  -1.79%  ' }
    @{ Cell = 'D31'; Value = '4.515' }
    @{ Cell = 'E31'; Value = '  -1.96%  ' }
    @{ Cell = 'D32'; Value = '1.524' }
    @{ Cell = 'E32'; Value = '  -0.60%  ' }
    @{ Cell = 'D33'; Value = '4.174' }
    @{ Cell = 'E33'; Value = '  -2.56%  ' }
    @{ Cell = 'D34'; Value = '0.05313' }
    @{ Cell = 'E34'; Value = '  -0.41%  ' }
    @{ Cell = 'D35'; Value = '1.229' }
    @{ Cell = 'E35'; Value = '  -1.14%  ' }
    @{ Cell = 'D36'; Value = '0.7380' }
    @{ Cell = 'E36'; Value = '  -2.39%  ' }
    @{ Cell = 'D37'; Value = '0.9988' }
    @{ Cell = 'E37'; Value = '  +0.03%  ' }
    @{ Cell = 'D38'; Value = '2.703' }
    @{ Cell = 'E38'; Value = '  +0.43%  ' }
    @{ Cell = 'D39'; Value = '0.01935' }
    @{ Cell = 'E39'; Value = '  -1.44%  ' }
    @{ Cell = 'D40'; Value = '2.735' }
    @{ Cell = 'E40'; Value = '  -0.59%  ' }
    @{ Cell = 'D41'; Value = '0.4474' }
    @{ Cell = 'E41'; Value = '  -0.17%  ' }
    @{ Cell = 'D42'; Value = '0.8758' }
    @{ Cell = 'E42'; Value = '  +1.99%  ' }
    @{ Cell = 'E43'; Value = '  -2.50%  ' }
    @{ Cell = 'D44'; Value = '71.16' }
    @{ Cell = 'E44'; Value = '  -1.85%  ' }
    @{ Cell = 'D45'; Value = '1.043.24' }
    @{ Cell = 'E45'; Value = '  -6.24%  ' }
    @{ Cell = 'E46'; Value = '  -0.09%  ' }
    @{ Cell = 'D47'; Value = '103.81' }
    @{ Cell = 'E47'; Value = '  +0.67%  ' }
    @{ Cell = 'D48'; Value = '1.819' }
    @{ Cell = 'E48'; Value = '  -1.11%  ' }
    @{ Cell = 'D49'; Value = '7.455' }
    @{ Cell = 'E49'; Value = '  -2.96%  ' }
    @{ Cell = 'D50'; Value = '9.551' }
    @{ Cell = 'E50'; Value = '  -0.09%  ' }
    @{ Cell = 'D51'; Value = '2.015.27' }
    @{ Cell = 'E51'; Value = '  -0.39%  ' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.NumberFormat = "General"
}
